$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '80.551.04'
$ws.Range('E2').Value = '  +5.31%  '
$ws.Range('D3').Value = '3.166.44'
$ws.Range('E3').Value = '  +2.59%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '208.92'
$ws.Range('E5').Value = '  +4.92%  '
$ws.Range('D6').Value = '623.31'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').Value = '0.271'
$ws.Range('E7').Value = '  +25.51%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  +6.20%  '
$ws.Range('D10').Value = '3.174.28'
$ws.Range('E10').Value = '  +2.96%  '
$ws.Range('D11').Value = '0.582'
$ws.Range('E11').Value = '  +26.13%  '
$ws.Range('E12').Value = '  +26.28%  '
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').Value = '3.756.80'
$ws.Range('E14').Value = '  +2.87%  '
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('D16').Value = '31.67'
$ws.Range('E16').Value = '  +6.90%  '
$ws.Range('D17').Value = '80.518.91'
$ws.Range('E17').Value = '  +5.41%  '
$ws.Range('D18').Value = '3.179.04'
$ws.Range('D19').Value = '14.15'
$ws.Range('E19').Value = '  +4.19%  '
$ws.Range('E20').Value = '  +9.46%  '
$ws.Range('D21').Value = '9.13'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').Value = '435.70'
$ws.Range('D23').Value = '5.14'
$ws.Range('E23').Value = '  +13.22%  '
$ws.Range('D24').Value = '6.94'
$ws.Range('E24').Value = '  +7.46%  '
$ws.Range('D26').Value = '76.00'
$ws.Range('E26').Value = '  +4.36%  '
$ws.Range('D27').Value = '4.65'
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('D28').Value = '10.79'
$ws.Range('E28').Value = '  +4.23%  '
$ws.Range('D29').Value = '0.996'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('D30').Value = '0.0000120'
$ws.Range('E30').Value = '  +7.64%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').Value = '8.85'
$ws.Range('E32').Value = '  +5.38%  '
$ws.Range('D33').Value = '556.23'
$ws.Range('E33').Value = '  +9.15%  '
$ws.Range('E34').Value = '  +0.50%  '
$ws.Range('D35').Value = '0.150'
$ws.Range('E35').Value = '  +13.04%  '
$ws.Range('E36').Value = '  +2.06%  '
$ws.Range('D37').Value = '22.84'
$ws.Range('E37').Value = '  +9.13%  '
$ws.Range('E38').Value = '  +19.76%  '
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('D40').Value = '0.402'
$ws.Range('E40').Value = '  +6.02%  '
$ws.Range('D41').Value = '20.79'
$ws.Range('E41').Value = '  +3.57%  '
$ws.Range('D42').Value = '163.45'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D44').Value = '5.59'
$ws.Range('E44').Value = '  +5.93%  '
$ws.Range('D45').Value = '189.37'
$ws.Range('E45').Value = '  -3.07%  '
$ws.Range('D46').Value = '1.79'
$ws.Range('E46').Value = '  +6.87%  '
$ws.Range('E47').Value = '  +7.25%  '
$ws.Range('D48').Value = '0.776'
$ws.Range('E48').Value = '  -2.58%  '
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').Value = '42.80'
$ws.Range('E50').Value = '  +3.97%  '
$ws.Range('D51').Value = '4.22'
$ws.Range('E51').Value = '  +6.75%  '
